$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 102, shifting existing rows 102:111 down to 103:112
$ws.Rows.Item(102).Insert()

# Populate the newly inserted row 102 with this week's data
$ws.Range("A102").Value = 10
$ws.Range("B102").Value = "Vega Modelo de Temuco"
$ws.Range("C102").Value = "La Araucanía"
$ws.Range("D102").Value = 44461
$ws.Range("E102").Value = 9
$ws.Range("F102").Value = 100114007
$ws.Range("G102").Value = "Jengibre"
$ws.Range("H102").Value = "Sin especificar"
$ws.Range("I102").Value = "Primera"
$ws.Range("J102").Value = 30
$ws.Range("K102").Value = 20000
$ws.Range("L102").Value = 20000
$ws.Range("M102").Value = 20000
$ws.Range("N102").Value = '$/caja 13 kilos'
$ws.Range("O102").Value = "Perú"
$ws.Range("P102").Value = 1538
$ws.Range("Q102").Value = 13
$ws.Range("R102").Value = "Hortaliza"

Write-Output "done"
